# Se agrega la información del tercer sensor de humedad.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sensores de humedad")

# Fill in the new MTH100 sensor row (row 7) with its characteristics
$ws.Range("C7").Value = "T: +/- 0.3°C;  RH: 3%RH"
$ws.Range("D7").Value = "t: <0.04°C/year; RH: <0.05 RH/year"
$ws.Range("E7").Value = "t: -20 to 80°C"
$ws.Range("F7").Value = "T: 0.1°C,    RH: 0.1%RH"
$ws.Range("G7").Value = "< 500 ohm (4-20 mA)"
$ws.Range("H7").Value = "15mA"
$ws.Range("I7").Value = "-"
$ws.Range("J7").Value = "ABS"
$ws.Range("K7").Value = "-"
$ws.Range("L7").Value = "T: 6t(63%): min = 2s,  max=5s"

# Fix a small typo in the existing MTH500H row (row 6): "H(90%)" -> "RH(90%)"
$ws.Range("L6").Value = "T: 6t(63%): min = 5s,  max=30s;   RH(90%): 8s"

$ws.Range("M7").Value = "<500 ohm"
$ws.Range("N7").Value = "-"
$ws.Range("O7").Value = "https://srcsl.com/catalogoPDFs/SensoresProce/SensoresHumedad/MTH100.pdf"

# Adjust column widths to better fit the new content
$ws.Columns.Item(6).ColumnWidth = 17.41796875
$ws.Columns.Item(12).ColumnWidth = 32.251302083333336

# Update the sheet view (scroll position / zoom / selection) to match the edited state
$ws.Application.ActiveWindow.Zoom = 69
$ws.Range("A8").Select()
